# "new rates as of 3/18"
# Updates the rate table: new pricing for the "Morris, IL" and "Belleville, MI"
# rows, three new destination rows (Monroe Township NJ / Monrovia MD /
# Owatonna MN), a new column (15) of rate tiers, and the "Belleville, MI"
# destination label gaining an extra space ("Belleville,  MI").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Make sure every data cell in B2:P15 carries the workbook's numeric
# style (xf index 1 -> numFmtId 2, "0.00") the way the rest of the sheet
# already does - this also guarantees "blank but styled" cells show up
# for columns that stay empty for a given row.
# ---------------------------------------------------------------------
$ws.Range("B2:P15").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# Row 1 - header / tier numbers. Column P (tier 15) is new.
# ---------------------------------------------------------------------
$ws.Range("P1").Value = 15

# ---------------------------------------------------------------------
# Row 3 - "Morris, IL" - full repricing, plus new N3/O3/P3 tiers.
# ---------------------------------------------------------------------
$ws.Range("B3").Value = 290
$ws.Range("C3").Value = 570
$ws.Range("D3").Value = 850
$ws.Range("E3").Value = 1120
$ws.Range("F3").Value = 1400
$ws.Range("G3").Value = 1680
$ws.Range("H3").Value = 1945
$ws.Range("I3").Value = 2200
$ws.Range("J3").Value = 2475
$ws.Range("K3").Value = 2750
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 3240
$ws.Range("N3").Value = 3510
$ws.Range("O3").Value = 3780
$ws.Range("P3").Value = 4015

# ---------------------------------------------------------------------
# New rows 13-15 (added before the Belleville relabel below, so the new
# shared strings land in the same order the source workbook has them).
# ---------------------------------------------------------------------
$ws.Range("A13").Value = "Monroe Township, NJ"
$ws.Range("B13").Value = 470
$ws.Range("C13").Value = 880
$ws.Range("D13").Value = 1245
$ws.Range("E13").Value = 1580
$ws.Range("F13").Value = 1925
$ws.Range("G13").Value = 2280
$ws.Range("H13").Value = 2625
$ws.Range("I13").Value = 2920
$ws.Range("J13").Value = 3240
$ws.Range("K13").Value = 3550

$ws.Range("A14").Value = "Monrovia, MD"
$ws.Range("B14").Value = 470
$ws.Range("C14").Value = 880
$ws.Range("D14").Value = 1245
$ws.Range("E14").Value = 1580
$ws.Range("F14").Value = 1925
$ws.Range("G14").Value = 2280
$ws.Range("H14").Value = 2625
$ws.Range("I14").Value = 2920
$ws.Range("J14").Value = 3240
$ws.Range("K14").Value = 3550

# ---------------------------------------------------------------------
# Row 6 - destination label gains a second space, and is fully
# repriced, plus new N6/O6/P6 tiers.
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "Belleville,  MI"
$ws.Range("B6").Value = 465
$ws.Range("C6").Value = 800
$ws.Range("D6").Value = 1185
$ws.Range("E6").Value = 1540
$ws.Range("F6").Value = 1825
$ws.Range("G6").Value = 2160
$ws.Range("H6").Value = 2415
$ws.Range("I6").Value = 2760
$ws.Range("J6").Value = 3060
$ws.Range("K6").Value = 3400
$ws.Range("L6").Value = 3740
$ws.Range("M6").Value = 4080
$ws.Range("N6").Value = 4355
$ws.Range("O6").Value = 4690
$ws.Range("P6").Value = 5025

$ws.Range("A15").Value = "Owatonna, MN"
$ws.Range("B15").Value = 525
$ws.Range("C15").Value = 1135
$ws.Range("D15").Value = 1475
$ws.Range("E15").Value = 1580
$ws.Range("F15").Value = 1840
$ws.Range("G15").Value = 2190
$ws.Range("H15").Value = 2345
$ws.Range("I15").Value = 2600
$ws.Range("J15").Value = 2880
$ws.Range("K15").Value = 3130

# ---------------------------------------------------------------------
# Column A widened to fit the new, longer destination names; it is the
# only single-width column now (columns C:P share one width, matching
# the rest of the table).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.6
$ws.Range("P1").EntireColumn.ColumnWidth = 6.7

# ---------------------------------------------------------------------
# Selection / active cell follows the newly-added rows, same as Excel
# leaves it after typing the last row of data.
# ---------------------------------------------------------------------
$ws.Range("A15:XFD16").Select()
